$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 25-26; existing rows 25-38 shift down to 27-40.
$ws.Rows("25:26").Insert()

# New row 25: Packham's Triumph, Segunda, $/caja 18 kilos granel, O'Higgins, 2022-09-05
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C25").Value = "Arica y Parinacota"
$ws.Range("D25").Value = 44809
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100104
$ws.Range("H25").Value = "Frutos de pepita"
$ws.Range("I25").Value = 100104005
$ws.Range("J25").Value = "Pera"
$ws.Range("K25").Value = "Packham's Triumph"
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 19000
$ws.Range("O25").Value = 20000
$ws.Range("P25").Value = 19500
$ws.Range("Q25").Value = "`$/caja 18 kilos granel"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 1083
$ws.Range("T25").Value = 18

# New row 26: Winter Nelis, Segunda, $/caja 18 kilos granel, O'Higgins, 2022-09-05
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44809
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100104
$ws.Range("H26").Value = "Frutos de pepita"
$ws.Range("I26").Value = 100104005
$ws.Range("J26").Value = "Pera"
$ws.Range("K26").Value = "Winter Nelis"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 270
$ws.Range("N26").Value = 19000
$ws.Range("O26").Value = 20000
$ws.Range("P26").Value = 19500
$ws.Range("Q26").Value = "`$/caja 18 kilos granel"
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 1083
$ws.Range("T26").Value = 18
